$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPAC")

# Insert two new columns before column D to make room for the newest two
# reporting quarters; this shifts the existing D:K data (and its styles)
# to F:M automatically.
$ws.Columns("D:E").Insert()

# The new D:E columns come in with default/general formatting - copy the
# number formats from column F (which holds the just-shifted original
# column D formatting) across D:E so dates/numbers render correctly.
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$ws.Columns("D:E").ColumnWidth = 13.83

# --- Period-ending date headers (rows 7, 38, 80) ---
$dateCols = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
foreach ($r in @(7, 38, 80)) {
    for ($i = 0; $i -lt $dateCols.Length; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $dateCols[$i]
    }
}

# --- Data rows: columns D (4) through M (13) ---
$rowsData = @{
    8 = @(102700, 96100, 86700, 95100, 98200, 100600, 84700, 84800, 95300, 95000)
    9 = @(65000, 61100, 55200, 58800, 58100, 60400, 52100, 50600, 56100, 54400)
    10 = @(37700, 35000, 31500, 36300, 40100, 40100, 32700, 34200, 39300, 40700)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(0, "NA", "NA", "NA", 14300, "NA", "NA", "NA", "NA", "NA")
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(85300, 75300, 70900, 76200, 90900, 77900, 70500, 69200, 73900, 71200)
    18 = @(17400, 20800, 15800, 18900, 7300, 22700, 14300, 15700, 21400, 23800)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(-20300, -6300, -5400, -5900, -5100, -5600, -5200, -5200, -5800, -4400)
    21 = @(7200, 24300, 20100, 22500, 12400, 26700, 17500, 19700, 24700, 28300)
    22 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    23 = @(-2900, 14500, 10400, 13000, 2200, 17100, 9000, 10400, 15600, 19300)
    24 = @(600, 4400, 3300, 4000, 2500, 5500, 2600, 3600, 10400, 5500)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(-3500, 10100, 7100, 9000, -300, 11600, 6400, 6800, 5200, 13900)
    27 = @(-3500, 10100, 7500, 9000, 3400, 11700, 6500, 7000, 5200, 14200)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0, 0, 0, -200, -2000, -700)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(20300, 6300, 5400, 5900, 5100, 5600, 5200, 5200, 5800, 4400)
    33 = @(-3500, 10100, 7500, 9000, 3400, 11700, 6500, 6800, 3200, 13500)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(-3500, 10100, 7500, 9000, 3400, 11700, 6500, 6800, 3200, 13500)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(14800, 25800, 20000, 9900, 14800, 36700, 14900, 5600, 24300, 52900)
    42 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    43 = @(41200, 43600, 41500, 38200, 38400, 41800, 42400, 36700, 38700, 42100)
    44 = @(128100, 125000, 119200, 108400, 112500, 102500, 103000, 104700, 105000, 107300)
    45 = @(1700, 6100, 5200, 3900, 1200, 4900, 4600, 4500, 104900, 4600)
    46 = @(185800, 200500, 186000, 160400, 166900, 185900, 164900, 151500, 272800, 206900)
    47 = @(13200, 15800, 13600, 12800, 11400, 16600, 26200, 28000, 29000, 45900)
    48 = @(649200, 651400, 656200, 660500, 666000, 670300, 671400, 679900, 688500, 756000)
    49 = @(14200, 4600, 3900, 4000, 4000, 13500, 13600, 13500, 13000, 26300)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(1100, 300, 200, 200, 300, 2400, 2400, 2500, 2500, 7800)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(863400, 872700, 859800, 838000, 848600, 888700, 878500, 875500, 1005800, 1042900)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(21000, 18800, 19800, 15100, 53700, 21200, 13300, 15500, 19400, 16200)
    58 = @(18300, "NA", "NA", "NA", "NA", "NA", "NA", 0, 0, 0)
    59 = @(38700, 85300, 38200, 26100, 8100, 30300, 31300, 28400, 35300, 83800)
    60 = @(78000, 104200, 57900, 41300, 61800, 51500, 44700, 43900, 54700, 100000)
    61 = @(308400, 296600, 294000, 292800, 291100, 292900, 291700, 292400, 302300, 306100)
    62 = @(39400, 40100, 37200, 40400, 41300, 45700, 47900, 47500, 49000, 40300)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(425800, 440800, 389000, 374400, 394300, 393900, 388100, 387400, 440200, 480600)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(153000, 146800, 186100, 179600, 171300, 214100, 210900, 207900, 200100, 197000)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(437700, 431800, 470800, 463500, 454300, 494800, 490400, 488100, 565600, 562300)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    81 = @(-3500, 10100, 7500, 9000, 3400, 11700, 6500, 6800, 3200, 13500)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(10000, 9800, 9700, 9600, 10200, 9600, 8400, 9300, 9200, 9000)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(20600, 16500, 12100, 12200, 11700, 34000, 10100, 13100, 25700, 18900)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(-6100, -8700, -2700, -6700, -7600, -7800, -2100, -6100, -7800, -1800)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-14800, -6700, -2500, -5800, 4900, -10000, 900, -17300, -7500, -1600)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(-44500, 0, 0, -7300, -37600, 0, 0, -100, -46600, 0)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(-16900, -4000, 500, -11300, -37600, -4000, 100, -14500, -46600, -3000)
    101 = @(100, 0, 0, 0, 0, 0, 0, 0, -200, 900)
    102 = @(-11000, 5800, 10100, -4900, -21900, 21800, 9300, -18700, -28600, 15100)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item($r, 4 + $i).Value = $v
        }
    }
}
